# Updated cryptos list with latest price/volume data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.918.23'
$ws.Range("E2").Value = '  +3.22%  '
$ws.Range("D3").Value = '1.882.83'
$ws.Range("E3").Value = '  +3.39%  '
$ws.Range("E4").Value = '  +0.13%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '327.84'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.31%  '
$ws.Range("E6").Value = '  +0.10%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4669'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +1.26%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3963'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +3.15%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07939'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +1.66%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9785'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +2.26%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '22.47'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +3.08%  '
$ws.Range("D12").Value = '1.845.27'
$ws.Range("E12").Value = '  +6.38%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.773'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +2.04%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.024'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +2.49%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.06987'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +1.80%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '89.00'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +3.17%  '
$ws.Range("E17").Value = '  +0.04%  '
$ws.Range("E18").Value = '  +2.44%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '17.06'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +1.82%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.004'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +0.31%  '
$ws.Range("D21").Value = '28.911.13'
$ws.Range("E21").Value = '  +3.18%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.369'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +1.07%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.15'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +1.92%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.116'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -0.34%  '
$ws.Range("D25").Value = '2.040.26'
$ws.Range("E25").Value = '  +3.39%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '153.65'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +1.31%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '19.49'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +1.80%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '5.795'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +1.86%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.016'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +2.40%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '120.08'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +3.26%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.09405'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +1.53%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.9462'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +1.14%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.339'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +1.52%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.356'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +3.97%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.352'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -2.12%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.05939'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -0.36%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02132'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -0.52%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.151'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +0.44%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '7.962'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +5.75%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.5756'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +3.11%  '
$ws.Range("B41").Value = 'Aptos'
$ws.Range("C41").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '10.03'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +1.40%  '
$ws.Range("B42").Value = 'Algorand'
$ws.Range("C42").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1801'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +1.99%  '
$ws.Range("B43").Value = 'EnergySwap'
$ws.Range("C43").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '11.94'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +3.14%  '
$ws.Range("B44").Value = 'Cronos'
$ws.Range("C44").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.07256'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +3.74%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.5363'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +2.37%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.151'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -3.27%  '
$ws.Range("E47").Value = '  -6.60%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.861'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +2.05%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '114.35'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +1.86%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.369'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +2.76%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.034'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +2.71%  '

Write-Output "Applied 103 cell updates"
